$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44487
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 2200
$ws.Range("L2").Value = 2200
$ws.Range("M2").Value = 2200
$ws.Range("P2").Value = 2200

$ws.Range("D3").Value = 44474
$ws.Range("J3").Value = 20
$ws.Range("K3").Value = 1600
$ws.Range("L3").Value = 1600
$ws.Range("M3").Value = 1600
$ws.Range("P3").Value = 1600

$ws.Range("D4").Value = 44497
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 2200
$ws.Range("L4").Value = 2200
$ws.Range("M4").Value = 2200
$ws.Range("P4").Value = 2200

$ws.Range("D5").Value = 44720
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 3600
$ws.Range("L5").Value = 3600
$ws.Range("M5").Value = 3600
$ws.Range("P5").Value = 3600

$ws.Range("D6").Value = 44679
$ws.Range("J6").Value = 30
$ws.Range("K6").Value = 5500
$ws.Range("L6").Value = 5500
$ws.Range("M6").Value = 5500
$ws.Range("P6").Value = 5500

$ws.Range("D7").Value = 44685
$ws.Range("J7").Value = 60
$ws.Range("K7").Value = 5000
$ws.Range("L7").Value = 6000
$ws.Range("M7").Value = 5333
$ws.Range("P7").Value = 5333

$ws.Range("D8").Value = 44719
$ws.Range("J8").Value = 80
$ws.Range("K8").Value = 3600
$ws.Range("L8").Value = 3600
$ws.Range("M8").Value = 3600
$ws.Range("P8").Value = 3600

$ws.Range("D9").Value = 44483
$ws.Range("J9").Value = 50
$ws.Range("K9").Value = 2200
$ws.Range("L9").Value = 2200
$ws.Range("M9").Value = 2200
$ws.Range("P9").Value = 2200

$ws.Range("D10").Value = 44484
$ws.Range("J10").Value = 40
$ws.Range("K10").Value = 2200
$ws.Range("L10").Value = 2200
$ws.Range("M10").Value = 2200
$ws.Range("P10").Value = 2200

$ws.Range("D11").Value = 44707
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 4700
$ws.Range("L11").Value = 4700
$ws.Range("M11").Value = 4700
$ws.Range("P11").Value = 4700

$ws.Range("D12").Value = 44677
$ws.Range("J12").Value = 20
$ws.Range("K12").Value = 5500
$ws.Range("L12").Value = 5500
$ws.Range("M12").Value = 5500
$ws.Range("P12").Value = 5500

$ws.Range("D13").Value = 44447
$ws.Range("J13").Value = 75
$ws.Range("K13").Value = 2200
$ws.Range("L13").Value = 2200
$ws.Range("M13").Value = 2200
$ws.Range("P13").Value = 2200

$ws.Range("D14").Value = 44453
$ws.Range("J14").Value = 20
$ws.Range("K14").Value = 2300
$ws.Range("L14").Value = 2300
$ws.Range("M14").Value = 2300
$ws.Range("P14").Value = 2300

$ws.Range("D15").Value = 44476
$ws.Range("J15").Value = 30
$ws.Range("K15").Value = 2200
$ws.Range("L15").Value = 2200
$ws.Range("M15").Value = 2200
$ws.Range("P15").Value = 2200

$ws.Range("D16").Value = 44669
$ws.Range("J16").Value = 60
$ws.Range("K16").Value = 6250
$ws.Range("L16").Value = 6250
$ws.Range("M16").Value = 6250
$ws.Range("P16").Value = 6250

$ws.Range("D17").Value = 44452
$ws.Range("J17").Value = 120
$ws.Range("K17").Value = 2300
$ws.Range("L17").Value = 2300
$ws.Range("M17").Value = 2300
$ws.Range("P17").Value = 2300

$ws.Range("D18").Value = 44706
$ws.Range("J18").Value = 90
$ws.Range("K18").Value = 4700
$ws.Range("L18").Value = 4700
$ws.Range("M18").Value = 4700
$ws.Range("P18").Value = 4700

$ws.Range("D19").Value = 44473
$ws.Range("J19").Value = 140
$ws.Range("K19").Value = 1600
$ws.Range("L19").Value = 1600
$ws.Range("M19").Value = 1600
$ws.Range("P19").Value = 1600

$ws.Range("D20").Value = 44203
$ws.Range("J20").Value = 30
$ws.Range("K20").Value = 2000
$ws.Range("L20").Value = 2000
$ws.Range("M20").Value = 2000
$ws.Range("P20").Value = 2000

$ws.Range("D21").Value = 44496
$ws.Range("J21").Value = 40
$ws.Range("K21").Value = 2200
$ws.Range("L21").Value = 2200
$ws.Range("M21").Value = 2200
$ws.Range("P21").Value = 2200

